$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650996137553771"
$ws1.Range("B2").Value = "go_stims-1650996137513734.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961375377345.csv"
$ws1.Range("B4").Value = "go_stims-16509961375377345.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996137553771.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509961396017907"
$ws2.Range("B2").Value = "ZB-match_5-16509961378497498.csv"
$ws2.Range("B3").Value = "TB-16509961391857285.csv"
$ws2.Range("B4").Value = "ZB-match_5-16509961382177641.csv"
$ws2.Range("B5").Value = "TB-16509961393697736.csv"
$ws2.Range("B6").Value = "OB-1650996138897732.csv"
$ws2.Range("B7").Value = "TB-16509961395777712.csv"
$ws2.Range("B8").Value = "ZB-match_0-16509961384977622.csv"
$ws2.Range("B9").Value = "OB-16509961391457675.csv"
$ws2.Range("B10").Value = "OB-16509961385697694.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509961396017907"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961396497655"
$ws4.Range("B2").Value = "MM_stims-16509961396177287.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961396017907.csv"
$ws4.Range("B4").Value = "MM_stims-1650996139633766.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961396177287.csv"
$ws4.Range("B6").Value = "MM_stims-16509961396497655.csv"
$ws4.Range("B7").Value = "ZM_stims-1650996139633766.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509961397137408"
$ws5.Range("B2").Value = "SAT_stims-16509961396497655.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961396817706.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961396977658.csv"
$ws5.Range("B5").Value = "SAT_stims-165099613966573.csv"
